# ToLTE.xlsx - "prepare import file for G2L nbr configuration"
#
# The template's second column header changes from the placeholder
# "target cell" to "target earfcn", the column is widened so the new,
# longer header still fits, and the sheet's last-saved cell selection
# moves on one row (C6 -> C7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header text: "target cell" -> "target earfcn"
$ws.Range("B1").Value = "target earfcn"

# Column B needs to be a bit wider to fit the new header text.
$ws.Columns.Item(2).ColumnWidth = 13

# The saved workbook's active cell moves from C6 down to C7.
$ws.Range("C7").Select()
